$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared string rich-text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "48"
$ws.Range("C9").Characters(27, 10).Text = "11/24/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/30/2025"

# --- Cells changing from numeric to text "0" (style 14 -> style 13, shared string) ---
# Use a leading apostrophe to force text, reusing the existing "0" shared string,
# then copy number-format/style from an existing style-13 "0" text cell (D14) so the
# resulting cell style matches (General format, no quote-prefix).
$ws.Range("C14").Value = "'0"
$ws.Range("F27").Value = "'0"
$ws.Range("C29").Value = "'0"
$ws.Range("C30").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)

# --- Cells changing from text to numeric (style 13 -> style 14 / 15) ---
$ws.Range("D16").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = 50
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# --- Remaining simple numeric value updates ---
$ws.Range("D15").Value = 1
$ws.Range("J15").Value = 27
$ws.Range("K15").Value = -25.925925925925
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = -43.274853801169
$ws.Range("L16").Value = -32.638888888888
$ws.Range("M16").Value = -57.079646017699
$ws.Range("N16").Value = -89.958592132505
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -15.789473684210
$ws.Range("I17").Value = 277
$ws.Range("J17").Value = 298
$ws.Range("K17").Value = -7.046979865771
$ws.Range("L17").Value = 5.725190839694
$ws.Range("M17").Value = 56.497175141242
$ws.Range("N17").Value = -43.121149897330
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 155
$ws.Range("K18").Value = -37.419354838709
$ws.Range("L18").Value = -19.166666666666
$ws.Range("M18").Value = -57.826086956521
$ws.Range("N18").Value = -90.934579439252
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 439
$ws.Range("K19").Value = -26.879271070615
$ws.Range("L19").Value = -38.857142857142
$ws.Range("M19").Value = -6.413994169096
$ws.Range("N19").Value = -24.824355971897
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -27.272727272727
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = -30.714285714285
$ws.Range("L20").Value = -41.566265060241
$ws.Range("M20").Value = -11.009174311926
$ws.Range("N20").Value = -89.281767955801
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -31.578947368421
$ws.Range("F21").Value = 68
$ws.Range("H21").Value = -15
$ws.Range("I21").Value = 911
$ws.Range("J21").Value = 1233
$ws.Range("K21").Value = -26.115166261151
$ws.Range("L21").Value = -26.234817813765
$ws.Range("M21").Value = -17.927927927927
$ws.Range("N21").Value = -76.688843398157
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 23.076923076923
$ws.Range("M22").Value = 14.285714285714
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -62.5
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = -4.761904761904
$ws.Range("I24").Value = 662
$ws.Range("J24").Value = 876
$ws.Range("K24").Value = -24.429223744292
$ws.Range("L24").Value = -41.519434628975
$ws.Range("M24").Value = -10.054347826087
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -77.777777777777
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -61.764705882352
$ws.Range("I25").Value = 191
$ws.Range("J25").Value = 363
$ws.Range("K25").Value = -47.382920110192
$ws.Range("L25").Value = -69.921259842519
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 58
$ws.Range("H26").Value = -27.586206896551
$ws.Range("I26").Value = 498
$ws.Range("J26").Value = 528
$ws.Range("K26").Value = -5.681818181818
$ws.Range("L26").Value = 4.184100418410
$ws.Range("M26").Value = -13.541666666666
$ws.Range("D27").Value = 1
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -15.151515151515
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 52
$ws.Range("J28").Value = 55
$ws.Range("K28").Value = -5.454545454545
$ws.Range("L28").Value = -34.177215189873
$ws.Range("N29").Value = -85.714285714285
$ws.Range("N30").Value = -85.185185185185
